$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A ("Trimestre" - quarter number).
# This shifts: old A -> B, old B -> C, old C -> D, and adjusts formulas/col widths automatically.
$ws.Columns("A").Insert()

# Header for the new column, copy formatting (bold, border, center, top) from the
# neighboring header cell (now B1) so it matches the other header cells.
$ws.Range("A1").Value = "Trimestre"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the new column with the repeating quarter number (1,2,3,4) for every block of 4 rows.
for ($r = 2; $r -le 37; $r++) {
    $q = (($r - 2) % 4) + 1
    $ws.Cells.Item($r, 1).Value = $q
}

# The old "Lluvia" column (now column C) gets right alignment, at the column level,
# and the header (C1) gets right alignment too (in addition to its existing bold/border/top).
$ws.Columns("C").HorizontalAlignment = -4152  # xlRight

# New "Lluvia" (rain) values were manually entered for the forecast rows (34-37),
# replacing the old placeholder values that were just repeats of rows 30-33.
$ws.Range("C34").Value = 0.23297599999999999
$ws.Range("C35").Value = 0.28274199999999999
$ws.Range("C36").Value = 0.187975
$ws.Range("C37").Value = 0.36432599999999998

# Row 35's rain cell ends up with a slightly different alignment (right + vertically
# centered) than the rest of the column (right only).
$ws.Range("C35").VerticalAlignment = -4108  # xlCenter

# Update the view: scroll position and active selection as recorded in the saved file.
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C37").Select()
